$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1559.3043
$ws.Range("I28").Value = 1186.1765
$ws.Range("J28").Value = 2616.5
$ws.Range("K28").Value = 1186.1765
$ws.Range("L28").Value = 2616.5
$ws.Range("M28").Value = -701.1765
$ws.Range("N28").Value = -3586.5

$ws.Range("H32").Value = 1914.125
$ws.Range("J32").Value = 1718.8334
$ws.Range("L32").Value = 1718.8334
$ws.Range("N32").Value = -2370.8334

$ws.Range("H63").Value = 47203.25
$ws.Range("J63").Value = 47203.25
$ws.Range("L63").Value = 47203.25
$ws.Range("N63").Value = -48451.25

$ws.Range("H66").Value = 47203.25
$ws.Range("J66").Value = 47203.25
$ws.Range("L66").Value = 141609.75
$ws.Range("N66").Value = -147849.75

$ws.Range("H106").Value = 280917.22
$ws.Range("I106").Value = 419400.84
$ws.Range("J106").Value = 3950
$ws.Range("K106").Value = 419400.84
$ws.Range("L106").Value = 3950
$ws.Range("M106").Value = -418769.84
$ws.Range("N106").Value = -5212

$ws.Range("H113").Value = 2246.577
$ws.Range("I113").Value = 2326.5789
$ws.Range("J113").Value = 2029.4286
$ws.Range("K113").Value = 2326.5789
$ws.Range("L113").Value = 2029.4286
$ws.Range("M113").Value = 927.4211
$ws.Range("N113").Value = -8537.428599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10305.532
$ws.Range("I32").Value = 9921.694
$ws.Range("J32").Value = 11752.308
$ws.Range("K32").Value = 9921.694
$ws.Range("L32").Value = 11752.308
$ws.Range("M32").Value = -9634.694
$ws.Range("N32").Value = -12326.308

$ws.Range("H74").Value = 29415344
$ws.Range("I74").Value = 50001484
$ws.Range("K74").Value = 50001484
$ws.Range("M74").Value = -50000610

$ws.Range("H77").Value = 29415344
$ws.Range("I77").Value = 50001484
$ws.Range("K77").Value = 250007420
$ws.Range("M77").Value = -250003052

$ws.Range("H122").Value = 5731.0435
$ws.Range("I122").Value = 6918.4116
$ws.Range("J122").Value = 2366.8333
$ws.Range("K122").Value = 20755.2348
$ws.Range("L122").Value = 7100.499899999999
$ws.Range("M122").Value = -18305.2348
$ws.Range("N122").Value = -12000.4999

$ws.Range("H132").Value = 13160779
$ws.Range("I132").Value = 25002388
$ws.Range("J132").Value = 3434.6667
$ws.Range("K132").Value = 75007164
$ws.Range("L132").Value = 10304.0001
$ws.Range("M132").Value = -75004634
$ws.Range("N132").Value = -15364.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 43122.285
$ws.Range("J76").Value = 43122.285
$ws.Range("L76").Value = 43122.285
$ws.Range("N76").Value = -43752.285

$ws.Range("H79").Value = 43122.285
$ws.Range("J79").Value = 43122.285
$ws.Range("L79").Value = 43122.285
$ws.Range("N79").Value = -45306.285

$ws.Range("H134").Value = 3866.2144
$ws.Range("I134").Value = 2940.8064
$ws.Range("J134").Value = 6474.1816
$ws.Range("K134").Value = 8822.4192
$ws.Range("L134").Value = 19422.5448
$ws.Range("M134").Value = -6287.4192
$ws.Range("N134").Value = -24492.5448

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1088.625
$ws.Range("I16").Value = 1027.8667
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1027.8667
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -740.8667
$ws.Range("N16").Value = -2574

$ws.Range("H58").Value = 4857.143
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 7750
$ws.Range("K58").Value = 1000
$ws.Range("L58").Value = 7750
$ws.Range("M58").Value = -797
$ws.Range("N58").Value = -8156

$ws.Range("H113").Value = 1088.625
$ws.Range("I113").Value = 1027.8667
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1027.8667
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1142.1333
$ws.Range("N113").Value = -6340

$ws.Range("H122").Value = 2090.8572
$ws.Range("I122").Value = 2106
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6318
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3868
$ws.Range("N122").Value = -10900

$ws.Range("H134").Value = 1961.5
$ws.Range("I134").Value = 2104.125
$ws.Range("J134").Value = 1771.3334
$ws.Range("K134").Value = 6312.375
$ws.Range("L134").Value = 5314.0002
$ws.Range("M134").Value = -3777.375
$ws.Range("N134").Value = -10384.0002

$ws.Range("H136").Value = 4857.143
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 7750
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 23250
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -28350

$ws.Range("H140").Value = 41495.8
$ws.Range("J140").Value = 41495.8
$ws.Range("L140").Value = 41495.8
$ws.Range("N140").Value = -51855.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38.72222
$ws.Range("J12").Value = 45.53846
$ws.Range("L12").Value = 136.61538
$ws.Range("N12").Value = -482.61538

$ws.Range("H33").Value = 43578692
$ws.Range("I33").Value = 37281964
$ws.Range("J33").Value = 66666696
$ws.Range("K33").Value = 223691784
$ws.Range("L33").Value = 400000176
$ws.Range("M33").Value = -223691501
$ws.Range("N33").Value = -400000742

$ws.Range("H64").Value = 2676.75
$ws.Range("I64").Value = 950
$ws.Range("J64").Value = 3252.3333
$ws.Range("K64").Value = 2850
$ws.Range("L64").Value = 9756.999899999999
$ws.Range("M64").Value = -2580
$ws.Range("N64").Value = -10296.9999

$ws.Range("H67").Value = 2676.75
$ws.Range("I67").Value = 950
$ws.Range("J67").Value = 3252.3333
$ws.Range("K67").Value = 2850
$ws.Range("L67").Value = 9756.999899999999
$ws.Range("M67").Value = -1914
$ws.Range("N67").Value = -11628.9999

$ws.Range("H80").Value = 2419.5
$ws.Range("I80").Value = 2496
$ws.Range("J80").Value = 2400.375
$ws.Range("K80").Value = 7488
$ws.Range("L80").Value = 7201.125
$ws.Range("M80").Value = -6552
$ws.Range("N80").Value = -9073.125

$ws.Range("H83").Value = 2419.5
$ws.Range("I83").Value = 2496
$ws.Range("J83").Value = 2400.375
$ws.Range("K83").Value = 22464
$ws.Range("L83").Value = 21603.375
$ws.Range("M83").Value = -17784
$ws.Range("N83").Value = -30963.375

$ws.Range("H86").Value = 2071.4285
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 2071.4285
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H113").Value = 899.8148
$ws.Range("I113").Value = 524.7143
$ws.Range("J113").Value = 1303.7693
$ws.Range("K113").Value = 1574.1429
$ws.Range("L113").Value = 3911.3079
$ws.Range("M113").Value = 595.8571000000002
$ws.Range("N113").Value = -8251.3079

$ws.Range("H132").Value = 843.7222
$ws.Range("I132").Value = 537.5
$ws.Range("J132").Value = 1088.7
$ws.Range("K132").Value = 4837.5
$ws.Range("L132").Value = 9798.300000000001
$ws.Range("M132").Value = -2307.5
$ws.Range("N132").Value = -14858.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3470.6
$ws.Range("I126").Value = 2421.44
$ws.Range("J126").Value = 4782.05
$ws.Range("K126").Value = 7264.32
$ws.Range("L126").Value = 14346.15
$ws.Range("M126").Value = -4794.32
$ws.Range("N126").Value = -19286.15

$ws.Range("H132").Value = 4552.0435
$ws.Range("I132").Value = 3223.75
$ws.Range("J132").Value = 6001.091
$ws.Range("K132").Value = 9671.25
$ws.Range("L132").Value = 18003.273
$ws.Range("M132").Value = -7141.25
$ws.Range("N132").Value = -23063.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4934.647
$ws.Range("I7").Value = 5257.8667
$ws.Range("J7").Value = 4679.4736
$ws.Range("K7").Value = 5257.8667
$ws.Range("L7").Value = 4679.4736
$ws.Range("M7").Value = -5145.8667
$ws.Range("N7").Value = -4903.4736

$ws.Range("H22").Value = 1005.52
$ws.Range("I22").Value = 425.9091
$ws.Range("J22").Value = 1460.9286
$ws.Range("K22").Value = 425.9091
$ws.Range("L22").Value = 1460.9286
$ws.Range("M22").Value = -130.9091
$ws.Range("N22").Value = -2050.9286

$ws.Range("H27").Value = 1005.52
$ws.Range("I27").Value = 425.9091
$ws.Range("J27").Value = 1460.9286
$ws.Range("K27").Value = 425.9091
$ws.Range("L27").Value = 1460.9286
$ws.Range("M27").Value = -318.9091
$ws.Range("N27").Value = -1674.9286

$ws.Range("H55").Value = 411.35294
$ws.Range("I55").Value = 310.2857
$ws.Range("J55").Value = 482.1
$ws.Range("K55").Value = 310.2857
$ws.Range("L55").Value = 482.1
$ws.Range("M55").Value = -137.2857
$ws.Range("N55").Value = -828.1

$ws.Range("H126").Value = 4934.647
$ws.Range("I126").Value = 5257.8667
$ws.Range("J126").Value = 4679.4736
$ws.Range("K126").Value = 15773.6001
$ws.Range("L126").Value = 14038.4208
$ws.Range("M126").Value = -13303.6001
$ws.Range("N126").Value = -18978.4208

$ws.Range("H132").Value = 11560.929
$ws.Range("I132").Value = 7986.8667
$ws.Range("J132").Value = 15684.846
$ws.Range("K132").Value = 23960.6001
$ws.Range("L132").Value = 47054.538
$ws.Range("M132").Value = -21430.6001
$ws.Range("N132").Value = -52114.538

$ws.Range("H139").Value = 60649.445
$ws.Range("J139").Value = 60649.445
$ws.Range("L139").Value = 60649.445
$ws.Range("N139").Value = -70929.44500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1047.409
$ws.Range("I132").Value = 527.8570999999999
$ws.Range("J132").Value = 1956.625
$ws.Range("K132").Value = 1583.5713
$ws.Range("L132").Value = 5869.875
$ws.Range("M132").Value = 946.4287000000002
$ws.Range("N132").Value = -10929.875
